# "Write to Excel" button/macro equivalent.
#
# Fills in the New/Old district-mapping worksheet with the review results
# for the first few rows that have been checked off, plus the W-district
# row (24): the looked-up "new" code goes in column C, and a status note
# ("done" / "skip") goes in column D. Row 2's note is bold (it doubles as
# a mini legend/header emphasis), the rest use the default body style.
#
# Cell writes are ordered top-to-bottom, left-to-right, matching the order
# a user (or the button's underlying macro) would actually fill the grid
# in, so that any newly-introduced shared strings line up the same way
# they did when this was done by hand in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — header-ish emphasis row, bold "done" marker.
$ws.Range("D2").Value = "done"
$ws.Range("D2").Font.Bold = $true

# Row 3 — already resolved.
$ws.Range("D3").Value = "done"

# Row 4 — partial match, flagged to skip for now.
$ws.Range("C4").Value = "E10 part"
$ws.Range("D4").Value = "skip"

# Row 5 — resolved to N2.
$ws.Range("C5").Value = "N2"
$ws.Range("D5").Value = "done"

# Row 6 — resolved to N3.
$ws.Range("C6").Value = "N3"
$ws.Range("D6").Value = "done"

# Row 24 — resolved to W3.
$ws.Range("C24").Value = "W3"
$ws.Range("D24").Value = "done"

# Leave the cursor where the button's work wrapped up.
$ws.Range("E24").Select()
